# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the freshly generated gh-pages data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - F-column updates
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    "F4"  = 111
    "F5"  = 1239
    "F6"  = 17514
    "F7"  = 326
    "F8"  = 231
    "F10" = 6645
    "F11" = 674
    "F13" = 98
    "F17" = 161
    "F18" = 52
    "F23" = 253
    "F24" = 946
    "F26" = 5116
    "F28" = 57
    "F29" = 11793
    "F32" = 187
    "F34" = 3891
    "F36" = 86
}
foreach ($cell in $expoUpdates.Keys) {
    $wsExpo.Range($cell).Value = $expoUpdates[$cell]
}

# Sheet "全部类型" (all types) - F-column updates
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    "F4"  = 111
    "F5"  = 1239
    "F6"  = 17514
    "F7"  = 326
    "F8"  = 231
    "F10" = 6645
    "F11" = 674
    "F13" = 98
    "F17" = 161
    "F18" = 52
    "F24" = 946
    "F26" = 5116
    "F29" = 57
    "F30" = 11793
    "F33" = 187
    "F35" = 3891
    "F37" = 86
}
foreach ($cell in $allUpdates.Keys) {
    $wsAll.Range($cell).Value = $allUpdates[$cell]
}
